$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "other" (column G) notes for existing rows (2-32) ---
$ws.Range("G2").Value = '3 DCCO, 2 RAZO'
$ws.Range("G3").Value = '10 DCCO, 1 RAZO'
$ws.Range("G4").Value = '1 BAEA, 6 corvid'
$ws.Range("G5").Value = '1 ATPU'
$ws.Range("G6").Value = '4 DCCO, 1 ATPU'
$ws.Range("G7").Value = '3 corvid, 1 RAZO'
$ws.Range("G8").Value = '3 RAZO'
$ws.Range("G9").Value = '1 NOGA'
$ws.Range("G10").Value = '1 COTE'
$ws.Range("G11").Value = '1 BAEA, 5 corvid'
$ws.Range("G12").Value = '1 RAZO'
$ws.Range("G14").Value = '2 ATPU'
$ws.Range("G15").Value = '1 ATPU, 1 RAZO, 4 corvid'
$ws.Range("G16").Value = '2 ATPU, 6 DCCO, 1 RAZO'
$ws.Range("G17").Value = '3 RAZO, 2 ATPU'
$ws.Range("G18").Value = '3 ATPU, 2 RAZO'
$ws.Range("G19").Value = '1 ATPU'
$ws.Range("G21").Value = '3 NOGA'
$ws.Range("G23").Value = '6 ATPU, 2 DCCO'
$ws.Range("G24").Value = '4 ATPU'
$ws.Range("G27").Value = '57 WISP, 1 ATPU, 5 porpoise, 3 RAZO'
$ws.Range("G28").Value = '15 ATPU, 1 COMU, 7 RAZO'

# --- New data row 37 needs the same date style as the rest of column A ---
# (rows 33-36 already carry that style from the pre-existing placeholder rows)
$ws.Range("A32").Copy()
$ws.Range("A37").PasteSpecial(-4122)

# --- New data rows 33-37 (A:F) ---
$ws.Range("A33").Value = 44395
$ws.Range("B33").Value = 584
$ws.Range("C33").Value = 41
$ws.Range("D33").Value = 35
$ws.Range("E33").Value = 4
$ws.Range("F33").Value = 219

$ws.Range("A34").Value = 44399
$ws.Range("B34").Value = 599
$ws.Range("C34").Value = 33
$ws.Range("D34").Value = 72
$ws.Range("E34").Value = 9
$ws.Range("F34").Value = 146
$ws.Range("G34").Value = '26 ATPU, 6 NOGA'

$ws.Range("A35").Value = 44400
$ws.Range("B35").Value = 565
$ws.Range("C35").Value = 53
$ws.Range("D35").Value = 83
$ws.Range("E35").Value = 2
$ws.Range("F35").Value = 151
$ws.Range("G35").Value = '30 WISP'

$ws.Range("A36").Value = 44401
$ws.Range("B36").Value = 596
$ws.Range("C36").Value = 42
$ws.Range("D36").Value = 96
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 212

$ws.Range("A37").Value = 44402
$ws.Range("B37").Value = 518
$ws.Range("C37").Value = 52
$ws.Range("D37").Value = 54
$ws.Range("E37").Value = 8
$ws.Range("F37").Value = 112

# --- Update selection to match the saved workbook state ---
$ws.Range("B39").Select()
